# Updates to power sector for curtailment and moving CHP out of flexible
# resources; updates to fuel balancing priorities.
#
# "municipal solid waste" (row 15) should no longer inherit the
# "natural gas peaker" (row 11) flag via formula (=B11); it is hard-coded
# to 0 (not a peaker / not a flexibility provider) on both the
# "Is This Plant Type a Peaker" sheet and the
# "Does This Plant Type Provide Flexibility" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("BPaFF-BITPTaP", "BPaFF-BDTPTPF")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Activate()
    $ws.Range("B15").Value = 0
    $ws.Range("B16").Select()
}

# Leave the "About" sheet as the active tab, matching the original workbook.
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
